# Daily attendance processing - 2025-10-06 21:38:47
# Update "Recorded By" (column G) values on the "Session Analysis Results" sheet.
# The update simply reorders the comma-separated list of recorders for
# specific rows (no content added/removed, only order changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "backup@backdoor.com, system, System"
    3   = "dnasr281@gmail.com, System"
    6   = "dnasr281@gmail.com, System"
    11  = "dnasr281@gmail.com, System"
    12  = "dnasr281@gmail.com, System"
    13  = "dnasr281@gmail.com, System"
    14  = "dnasr281@gmail.com, System"
    15  = "dnasr281@gmail.com, System"
    29  = "backup@backdoor.com, system, System"
    30  = "dnasr281@gmail.com, System"
    33  = "dnasr281@gmail.com, System"
    38  = "dnasr281@gmail.com, System"
    39  = "dnasr281@gmail.com, System"
    40  = "dnasr281@gmail.com, System"
    41  = "dnasr281@gmail.com, System"
    42  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, system, System"
    57  = "dnasr281@gmail.com, System"
    60  = "dnasr281@gmail.com, System"
    65  = "dnasr281@gmail.com, System"
    66  = "dnasr281@gmail.com, System"
    67  = "dnasr281@gmail.com, System"
    68  = "dnasr281@gmail.com, System"
    69  = "dnasr281@gmail.com, System"
    89  = "dnasr281@gmail.com, System"
    90  = "admin@admin.com, dnasr281@gmail.com"
    93  = "dnasr281@gmail.com, System"
    115 = "dnasr281@gmail.com, System"
    116 = "admin@admin.com, dnasr281@gmail.com"
    119 = "dnasr281@gmail.com, System"
    141 = "dnasr281@gmail.com, System"
    142 = "admin@admin.com, dnasr281@gmail.com"
    145 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
